# Fix ordre virement (carte_sejour,cin) condition issue
#
# Replaces the two existing data rows (2-3) with corrected beneficiary
# data and appends the remaining beneficiary rows (4-7) plus a totals
# row (8), matching the regenerated "etat des virements" export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing text storage for
# numeric-looking strings (CIN/account numbers that must keep leading
# zeros / not be coerced to a Number by Excel's auto-detection).
function Set-TextCell($ws, $ref, [string]$text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($ref).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

function Test-NumericText([string]$s) {
    return $s -match '^[+-]?\d+(\.\d+)?$'
}

function Set-RowValues($ws, [int]$r, [object[]]$vals) {
    $cols = @("A","B","C","D","E","F","G","H","I","J","K")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ref = "$($cols[$i])$r"
        $v = $vals[$i]
        $isNumeric = ($v -is [int]) -or ($v -is [int32]) -or ($v -is [int64]) -or ($v -is [double]) -or ($v -is [decimal])
        if ($isNumeric) {
            $ws.Range($ref).Value = $v
        } elseif (Test-NumericText $v) {
            Set-TextCell $ws $ref $v
        } else {
            $ws.Range($ref).Value = $v
        }
    }
}

$rows = @{
    2 = @("AGENCE KHATABI","354646","54544646446464646464444464","AGENCE 1","BMCI","Logement de fonction","908/LF/DIRECTION REGIONALE SUD","mensuelle",16000,1400,14600)
    3 = @("CHARIJI ABDELLAH","BJ36877","00101211111292695000201732","AOURIR","BMCE","Point de vente","389/AOURIR","mensuelle",7000,700,6300)
    4 = @("ACHENGLI LAILA","J207703","00101211115087750001201090","AIT SOUSS","BP","Direction régionale","908/DIRECTION REGIONALE SUD","mensuelle",20000,3000,17000)
    5 = @("ACHENGLI LAILA","J207703","00101211115087750001201090","AIT SOUSS","BP","Siège","900/PATIO","mensuelle",4500,450,4050)
    6 = @("AGENCE KHATABI","354646","54544646446464646464444464","AGENCE 1","BMCI","Supervision","001/SUP SUD","mensuelle",2400,0,2400)
    7 = @("NACER YASSINE","L234567","78017098772736274634834384","TOUHAMI","ATTIJARI WAFA BANK","Point de vente","805/KOUTOUBIA","mensuelle",12000,1800,10200)
    8 = @(" "," "," "," "," "," "," "," ",61900,7350,54550)
}

foreach ($r in 2..8) {
    Set-RowValues $ws $r $rows[$r]
}
